$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits: add missing F5/G5 observations to the accuracy log ---
$ws.Range("F5").Value = -0.14577999999999999
$ws.Range("G5").Value = -0.45687

# --- Formula edits: switch the COUNTA() denominators to COUNT() in the RMSE cells ---
$ws.Range("D100").Formula = "=SQRT(SUMSQ(D2:D97)/COUNT(D2:D97))"
$ws.Range("F100").Formula = "=SQRT(SUMSQ(F2:F97)/COUNT(F2:F97))"

# --- Column D width: widen to fit the "accuracy_log_od_mean" header ---
$ws.Columns.Item(4).ColumnWidth = 21.02

# --- View state: scroll down near the bottom of the data and select D100 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 85
$win.ScrollColumn = 1
$ws.Range("D100").Select()
